# Journal de bord CPNVoiturage - mise a jour du JdB et de la planif
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix wording of an existing task description (shared string index 47)
$ws.Range("F27").Value = "Tests de css pour avoir des onglets responsive"

# 2) Fill in the previously empty rows 28-31 with new journal entries
$ws.Range("C28").Value = "Implémentation complète du gabarit"
$ws.Range("D28").Value = 44245
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = "Connecter le gabarit correctement et ajout correct de la vue au choix"

$ws.Range("C29").Value = "Création "
$ws.Range("D29").Value = 44245
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = "Création et gestion du css de la page du profil"

$ws.Range("F30").Value = "Javascript pour disable tous les champs désactivés en fonction des jours"
$ws.Range("C30").Value = "JS pour les choxi des jours"
$ws.Range("D30").Value = 44245
$ws.Range("E30").Value = 45

$ws.Range("F31").Value = "Création de la page de login et du css"
$ws.Range("C31").Value = "Création de la page de login"
$ws.Range("D31").Value = 44246
$ws.Range("E31").Value = 100

# 3) Update the sheet view (scroll position and selected cell)
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F25").Select()
